# Update cryptocurrency price (D) and 1h volume change (E) columns
# with freshly scraped values from coinranking.com, per GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '23.311.73'
$ws.Cells.Item(2, 5).Value = '  -1.82%  '

$ws.Cells.Item(3, 4).Value = '1.628.47'
$ws.Cells.Item(3, 5).Value = '  -1.87%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  +0.07%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '1.001'
$ws.Cells.Item(5, 5).Value = '  +0.08%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '297.97'

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.3765'
$ws.Cells.Item(7, 5).Value = '  -1.26%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '49.88'
$ws.Cells.Item(8, 5).Value = '  -2.68%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.3475'
$ws.Cells.Item(9, 5).Value = '  -3.97%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.08031'
$ws.Cells.Item(10, 5).Value = '  -1.96%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.198'
$ws.Cells.Item(11, 5).Value = '  -2.82%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.001'
$ws.Cells.Item(12, 5).Value = '  +0.09%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '21.86'
$ws.Cells.Item(13, 5).Value = '  -3.19%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.277'
$ws.Cells.Item(14, 5).Value = '  -2.93%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.206'
$ws.Cells.Item(15, 5).Value = '  -3.04%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.00001188'
$ws.Cells.Item(16, 5).Value = '  -3.17%  '

$ws.Cells.Item(17, 4).Value = '1.628.09'

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '94.62'
$ws.Cells.Item(18, 5).Value = '  -3.19%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06957'
$ws.Cells.Item(19, 5).Value = '  -0.77%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.600'
$ws.Cells.Item(20, 5).Value = '  -3.12%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '17.23'
$ws.Cells.Item(21, 5).Value = '  -2.34%  '

$ws.Cells.Item(22, 5).Value = '  +0.22%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '12.33'
$ws.Cells.Item(23, 5).Value = '  -3.74%  '

$ws.Cells.Item(24, 4).Value = '23.303.05'

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.418'
$ws.Cells.Item(25, 5).Value = '  -3.35%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.956'
$ws.Cells.Item(26, 5).Value = '  -1.77%  '

$ws.Cells.Item(27, 5).Value = '  -1.99%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '149.68'
$ws.Cells.Item(28, 5).Value = '  -2.32%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.154'
$ws.Cells.Item(29, 5).Value = '  -1.48%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '130.64'
$ws.Cells.Item(30, 5).Value = '  -2.68%  '

$ws.Cells.Item(31, 4).Value = '1.808.07'
$ws.Cells.Item(31, 5).Value = '  -2.00%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '6.708'
$ws.Cells.Item(32, 5).Value = '  -6.73%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.116'
$ws.Cells.Item(33, 5).Value = '  -5.76%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '11.14'
$ws.Cells.Item(34, 5).Value = '  -7.56%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.9781'
$ws.Cells.Item(35, 5).Value = '  -7.30%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02651'
$ws.Cells.Item(36, 5).Value = '  -5.98%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.08737'
$ws.Cells.Item(37, 5).Value = '  -0.76%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2408'
$ws.Cells.Item(38, 5).Value = '  -4.39%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.792'
$ws.Cells.Item(39, 5).Value = '  -5.32%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.06727'
$ws.Cells.Item(40, 5).Value = '  -4.07%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '12.64'
$ws.Cells.Item(41, 5).Value = '  -3.05%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.6776'
$ws.Cells.Item(42, 5).Value = '  -3.31%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.286'
$ws.Cells.Item(43, 5).Value = '  -3.71%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '15.42'
$ws.Cells.Item(44, 5).Value = '  -4.47%  '

$ws.Cells.Item(45, 5).Value = '  +0.13%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.6280'
$ws.Cells.Item(46, 5).Value = '  -3.69%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.225'
$ws.Cells.Item(47, 5).Value = '  -3.49%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.884'
$ws.Cells.Item(48, 5).Value = '  -2.01%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.07611'
$ws.Cells.Item(49, 5).Value = '  -3.97%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '125.98'
$ws.Cells.Item(50, 5).Value = '  -1.89%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.215'
$ws.Cells.Item(51, 5).Value = '  +1.72%  '
